# Auto-generated Excel COM-interop script to update the cryptocurrency price table
# Applies the diff: updates Price (D) and Volume(1h) (E) columns for most rows,
# and swaps two pairs of rows (30/31 and 48/49) whose coin/link/price/volume content changed places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as literal text,
# even when the text looks numeric (e.g. "603.95", "1.00", "0.0000191").
# This mirrors the source data, where prices are inline strings, not numbers.
# A leading apostrophe forces Excel to treat the entry as text; the cells original
# style is restored afterwards so no stray "quote prefix" style is left behind.
function Set-TextValue {
    param($Worksheet, $Row, $Col, $Text)
    $cell = $Worksheet.Cells.Item($Row, $Col)
    $origStyle = $cell.Style
    $cell.Value = "'" + $Text
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws 2 4 "65.354.98"
$ws.Cells.Item(2, 5).Value = "  +3.08%  "

# Row 3
Set-TextValue $ws 3 4 "2.646.90"
$ws.Cells.Item(3, 5).Value = "  +1.90%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.00%  "

# Row 5
Set-TextValue $ws 5 4 "603.95"
$ws.Cells.Item(5, 5).Value = "  +2.64%  "

# Row 6
Set-TextValue $ws 6 4 "156.56"
$ws.Cells.Item(6, 5).Value = "  +5.03%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.04%  "

# Row 8
Set-TextValue $ws 8 4 "0.589"
$ws.Cells.Item(8, 5).Value = "  +1.01%  "

# Row 9
Set-TextValue $ws 9 4 "0.121"
$ws.Cells.Item(9, 5).Value = "  +10.57%  "

# Row 10
Set-TextValue $ws 10 4 "0.409"
$ws.Cells.Item(10, 5).Value = "  +6.13%  "

# Row 11
Set-TextValue $ws 11 4 "5.80"
$ws.Cells.Item(11, 5).Value = "  +1.28%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +2.53%  "

# Row 13
Set-TextValue $ws 13 4 "29.42"
$ws.Cells.Item(13, 5).Value = "  +6.81%  "

# Row 14
Set-TextValue $ws 14 4 "0.0000191"
$ws.Cells.Item(14, 5).Value = "  +22.54%  "

# Row 15
Set-TextValue $ws 15 4 "3.123.10"
$ws.Cells.Item(15, 5).Value = "  +1.87%  "

# Row 16
Set-TextValue $ws 16 4 "65.183.61"
$ws.Cells.Item(16, 5).Value = "  +3.07%  "

# Row 17
Set-TextValue $ws 17 4 "2.658.34"
$ws.Cells.Item(17, 5).Value = "  +2.76%  "

# Row 18
Set-TextValue $ws 18 4 "12.72"
$ws.Cells.Item(18, 5).Value = "  +5.49%  "

# Row 19
Set-TextValue $ws 19 4 "4.92"
$ws.Cells.Item(19, 5).Value = "  +5.03%  "

# Row 20
Set-TextValue $ws 20 4 "360.40"
$ws.Cells.Item(20, 5).Value = "  +4.73%  "

# Row 21
Set-TextValue $ws 21 4 "7.38"
$ws.Cells.Item(21, 5).Value = "  +8.05%  "

# Row 23
Set-TextValue $ws 23 4 "69.19"
$ws.Cells.Item(23, 5).Value = "  +3.95%  "

# Row 24
Set-TextValue $ws 24 4 "1.70"
$ws.Cells.Item(24, 5).Value = "  -0.42%  "

# Row 25
Set-TextValue $ws 25 4 "9.44"
$ws.Cells.Item(25, 5).Value = "  +3.09%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +1.21%  "

# Row 27
Set-TextValue $ws 27 4 "8.28"
$ws.Cells.Item(27, 5).Value = "  +0.61%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +3.13%  "

# Row 29
Set-TextValue $ws 29 4 "0.0₃0976"
$ws.Cells.Item(29, 5).Value = "  +14.60%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Bittensor"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws 30 4 "547.82"
$ws.Cells.Item(30, 5).Value = "  -1.30%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws 31 4 "2.20"
$ws.Cells.Item(31, 5).Value = "  +8.63%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -0.09%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +2.74%  "

# Row 34
Set-TextValue $ws 34 4 "5.58"
$ws.Cells.Item(34, 5).Value = "  +6.09%  "

# Row 35
Set-TextValue $ws 35 4 "6.42"
$ws.Cells.Item(35, 5).Value = "  +6.24%  "

# Row 36
Set-TextValue $ws 36 4 "0.430"
$ws.Cells.Item(36, 5).Value = "  +4.47%  "

# Row 37
Set-TextValue $ws 37 4 "20.68"
$ws.Cells.Item(37, 5).Value = "  +6.71%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +5.56%  "

# Row 39
Set-TextValue $ws 39 4 "161.99"
$ws.Cells.Item(39, 5).Value = "  -2.00%  "

# Row 40
Set-TextValue $ws 40 4 "0.999"
$ws.Cells.Item(40, 5).Value = "  -0.06%  "

# Row 41
Set-TextValue $ws 41 4 "1.00"
$ws.Cells.Item(41, 5).Value = "  +0.08%  "

# Row 42
Set-TextValue $ws 42 4 "42.64"
$ws.Cells.Item(42, 5).Value = "  +7.93%  "

# Row 43
Set-TextValue $ws 43 4 "167.03"
$ws.Cells.Item(43, 5).Value = "  +1.12%  "

# Row 44
Set-TextValue $ws 44 4 "4.18"
$ws.Cells.Item(44, 5).Value = "  +4.46%  "

# Row 45
Set-TextValue $ws 45 4 "0.0619"
$ws.Cells.Item(45, 5).Value = "  +7.32%  "

# Row 46
Set-TextValue $ws 46 4 "23.35"
$ws.Cells.Item(46, 5).Value = "  +2.31%  "

# Row 47
Set-TextValue $ws 47 4 "2.28"
$ws.Cells.Item(47, 5).Value = "  +8.75%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "VeChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws 48 4 "0.0265"
$ws.Cells.Item(48, 5).Value = "  +7.27%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Mantle"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws 49 4 "0.656"
$ws.Cells.Item(49, 5).Value = "  +4.07%  "

# Row 50
Set-TextValue $ws 50 4 "0.0984"
$ws.Cells.Item(50, 5).Value = "  +2.74%  "

# Row 51
Set-TextValue $ws 51 4 "19.76"
$ws.Cells.Item(51, 5).Value = "  +4.04%  "

